$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STEPS")
$ws.Activate() | Out-Null

# Remember the width of the column that will end up just to the left
# of the newly inserted column, so the new column inherits it (this is
# what Excel does natively when inserting a column).
$leftWidth = $ws.Columns("F:F").ColumnWidth

# Insert a new column before the current column G ("TC_STEP_ACTION"),
# shifting everything from G onward one column to the right.
$ws.Columns("G:G").Insert()
$ws.Columns("G:G").ColumnWidth = $leftWidth

# New header cell for the inserted column.
$ws.Range("G1").Value = "TC_STEP_CALL_DATASET"

# Update the active selection on the sheet.
$ws.Range("G2").Select() | Out-Null
